$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 4.125059742542763
$ws.Range("D2").Value = 6.231446479256613
$ws.Range("H2").Value = 4.864952974133973
$ws.Range("J2").Value = 3.572898773933583
$ws.Range("L2").Value = 1.830940742846879

# Row 3
$ws.Range("B3").Value = 2.724580040526347
$ws.Range("D3").Value = 2.333436841424388
$ws.Range("H3").Value = 3.115723239628307

# Row 4
$ws.Range("B4").Value = 5.203216707823461
$ws.Range("D4").Value = 6.486603965511964
$ws.Range("F4").Value = 3.579079026444063
$ws.Range("J4").Value = 5.543967131514354
